# Generate Report for handoff
#
# The source file "629e121d-9986-4cbd-a3f7-576432f1d3cc.md" is being handed
# off again: its Status flips from "Handed back" to "Not yet handed off" on
# the Overview sheet (for both locales) and on each per-locale detail sheet,
# and the "Latest Handoff Datetime" on each locale sheet is refreshed to the
# new handoff timestamp. The previous handback info (file / datetime) is left
# untouched since nothing has been handed back yet for this new handoff.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 629e121d...md file; columns B (zh-cn) / C (de-de) hold status text.
$overview.Range("B3").Value = "Not yet handed off"
$overview.Range("C3").Value = "Not yet handed off"

# zh-cn detail sheet: row 3 is the 629e121d...md file.
$zhcn.Range("B3").Value = "Not yet handed off"
$zhcn.Range("D3").Value = "2016-01-08 15:25:50"

# de-de detail sheet: row 3 is the 629e121d...md file.
$dede.Range("B3").Value = "Not yet handed off"
$dede.Range("D3").Value = "2016-01-08 15:26:05"
